$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap match data (F:V) between row 111 and row 113 ---
# Row 111
$ws.Range("F111").Value = "Legionowo"
$ws.Range("G111").Value = 0
$ws.Range("H111").Value = "Zambrow"
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 2.09
$ws.Range("K111").Value = "27/10/2023 01:12"
$ws.Range("L111").Value = 2.32
$ws.Range("M111").Value = "28/10/2023 12:59"
$ws.Range("N111").Value = 3.5
$ws.Range("O111").Value = "27/10/2023 01:12"
$ws.Range("P111").Value = 3.54
$ws.Range("Q111").Value = "28/10/2023 12:50"
$ws.Range("R111").Value = 2.67
$ws.Range("S111").Value = "27/10/2023 01:12"
$ws.Range("T111").Value = 2.55
$ws.Range("U111").Value = "28/10/2023 12:59"
$ws.Range("V111").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-i/legionowo-olimpia-zambrow/A52kGPr3/"

# Row 113
$ws.Range("F113").Value = "LKS Lomza"
$ws.Range("G113").Value = 2
$ws.Range("H113").Value = "Concordia Elblag"
$ws.Range("I113").Value = 4
$ws.Range("J113").Value = 2.08
$ws.Range("K113").Value = "27/10/2023 01:12"
$ws.Range("L113").Value = 2.2
$ws.Range("M113").Value = "28/10/2023 12:59"
$ws.Range("N113").Value = 3.33
$ws.Range("O113").Value = "27/10/2023 01:12"
$ws.Range("P113").Value = 3.54
$ws.Range("Q113").Value = "28/10/2023 12:59"
$ws.Range("R113").Value = 2.73
$ws.Range("S113").Value = "27/10/2023 01:12"
$ws.Range("T113").Value = 2.72
$ws.Range("U113").Value = "28/10/2023 12:59"
$ws.Range("V113").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-i/lks-lomza-concordia-elblag/MqMffrMq/"

# --- Cyclic rearrangement of match data (F:V) across rows 124-127 ---
# Row 124
$ws.Range("F124").Value = "Concordia Elblag"
$ws.Range("G124").Value = 4
$ws.Range("H124").Value = "Sulejowek"
$ws.Range("I124").Value = 1
$ws.Range("J124").Value = 2.67
$ws.Range("K124").Value = "03/11/2023 02:13"
$ws.Range("L124").Value = 2.81
$ws.Range("M124").Value = "04/11/2023 13:52"
$ws.Range("N124").Value = 3.23
$ws.Range("O124").Value = "03/11/2023 02:13"
$ws.Range("P124").Value = 3.35
$ws.Range("Q124").Value = "04/11/2023 13:52"
$ws.Range("R124").Value = 2.16
$ws.Range("S124").Value = "03/11/2023 02:13"
$ws.Range("T124").Value = 2.22
$ws.Range("U124").Value = "04/11/2023 13:52"
$ws.Range("V124").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-i/concordia-elblag-victoria-sulejowek/jPZt6uLe/"

# Row 125
$ws.Range("F125").Value = "GKS Belchatow"
$ws.Range("G125").Value = 2
$ws.Range("H125").Value = "T. Mazowiecki"
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 1.93
$ws.Range("K125").Value = "03/11/2023 02:13"
$ws.Range("L125").Value = 1.77
$ws.Range("M125").Value = "04/11/2023 13:41"
$ws.Range("N125").Value = 3.43
$ws.Range("O125").Value = "03/11/2023 02:13"
$ws.Range("P125").Value = 3.75
$ws.Range("Q125").Value = "04/11/2023 13:41"
$ws.Range("R125").Value = 2.96
$ws.Range("S125").Value = "03/11/2023 02:13"
$ws.Range("T125").Value = 3.66
$ws.Range("U125").Value = "04/11/2023 13:41"
$ws.Range("V125").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-i/gks-belchatow-tomaszow-mazowiecki/44ebE3DF/"

# Row 126
$ws.Range("F126").Value = "Mlawa"
$ws.Range("G126").Value = 2
$ws.Range("H126").Value = "LKS Lomza"
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 1.72
$ws.Range("K126").Value = "03/11/2023 02:13"
$ws.Range("L126").Value = 1.83
$ws.Range("M126").Value = "04/11/2023 13:52"
$ws.Range("N126").Value = 3.74
$ws.Range("O126").Value = "03/11/2023 02:13"
$ws.Range("P126").Value = 4.01
$ws.Range("Q126").Value = "04/11/2023 13:52"
$ws.Range("R126").Value = 3.33
$ws.Range("S126").Value = "03/11/2023 02:13"
$ws.Range("T126").Value = 3.25
$ws.Range("U126").Value = "04/11/2023 13:05"
$ws.Range("V126").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-i/mks-mlawa-lks-lomza/4Swo5az2/"

# Row 127
$ws.Range("F127").Value = "Pilica Bialobrzegi"
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = "Skierniewice"
$ws.Range("I127").Value = 4
$ws.Range("J127").Value = 2.88
$ws.Range("K127").Value = "03/11/2023 02:13"
$ws.Range("L127").Value = 3.2
$ws.Range("M127").Value = "03/11/2023 11:31"
$ws.Range("N127").Value = 3.39
$ws.Range("O127").Value = "03/11/2023 02:13"
$ws.Range("P127").Value = 3.56
$ws.Range("Q127").Value = "04/11/2023 12:02"
$ws.Range("R127").Value = 1.97
$ws.Range("S127").Value = "03/11/2023 02:13"
$ws.Range("T127").Value = 1.93
$ws.Range("U127").Value = "03/11/2023 11:31"
$ws.Range("V127").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-i/pilica-bialobrzegi-unia-skierniewice/pOGXS2S8/"

# --- Append new row 129 (copy formatting from row 128, then set values) ---
$ws.Range("A128:V128").Copy()
$ws.Range("A129").PasteSpecial(-4122)
$ws.Range("A129").Value = 128
$ws.Range("B129").Value = "poland"
$ws.Range("C129").Value = "iii-liga-group-i"
$ws.Range("D129").Value = "2023-2024"
$ws.Range("E129").Value = 45240.75
$ws.Range("F129").Value = "Sulejowek"
$ws.Range("G129").Value = 1
$ws.Range("H129").Value = "Mlawa"
$ws.Range("I129").Value = 1
$ws.Range("J129").Value = 1.69
$ws.Range("K129").Value = "10/11/2023 06:42"
$ws.Range("L129").Value = 1.67
$ws.Range("M129").Value = "10/11/2023 17:55"
$ws.Range("N129").Value = 3.8
$ws.Range("O129").Value = "10/11/2023 06:42"
$ws.Range("P129").Value = 4.09
$ws.Range("Q129").Value = "10/11/2023 17:55"
$ws.Range("R129").Value = 3.7
$ws.Range("S129").Value = "10/11/2023 06:42"
$ws.Range("T129").Value = 3.82
$ws.Range("U129").Value = "10/11/2023 17:55"
$ws.Range("V129").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-i/victoria-sulejowek-mks-mlawa/CjnD8Ytf/"
